# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "municipio-nombre" column (I) was previously annotated as a measure
# (iaest-measure:municipio-nombre / medida / xsd:int). With the newly
# curated dimensions it must instead be annotated as a dimension that
# references an area, matching the treatment already used for the
# "provincia-nombre" (J) and "comarca-nombre" (M) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: metadata "type" row -> dimension/measure identifier
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# Row 3: metadata "dim"/"medida" marker
$ws.Range("I3").Value = "dim"

# Row 4: metadata datatype/URI marker
$ws.Range("I4").Value = "URI-Municipio"
